# "reload 3 times when error"
#
# The three "login attempt" rows (row 2, row 3, row 4 -- usernames
# 100160016W, 100160423W, 100160403W) are updated with new visit counts
# and new "last visit" timestamps, reflecting that the scraper reloaded
# those accounts (up to) three more times.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 100160016W / Abcfinance*20
$ws.Range("C2").Value = 200
$ws.Range("D2").Value = "2020-03-13 00:18:39.043948"

# Row 3 - 100160423W / Idfjobs*0505
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = "2020-03-13 01:59:22.479527"

# Row 4 - 100160403W / Bestshore@05
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "2020-03-13 02:08:39.637561"
